$wb = $excel.ActiveWorkbook

# Updated "想去人数" (F column) values for sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1634
$ws1.Range("F5").Value = 629
$ws1.Range("F6").Value = 1102
$ws1.Range("F8").Value = 11582
$ws1.Range("F9").Value = 29
$ws1.Range("F12").Value = 372
$ws1.Range("F13").Value = 1094
$ws1.Range("F15").Value = 12410
$ws1.Range("F16").Value = 13151
$ws1.Range("F17").Value = 34
$ws1.Range("F21").Value = 244
$ws1.Range("F24").Value = 130

# Updated "想去人数" (F column) values for sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1634
$ws4.Range("F5").Value = 629
$ws4.Range("F6").Value = 1102
$ws4.Range("F8").Value = 11582
$ws4.Range("F9").Value = 29
$ws4.Range("F11").Value = 455
$ws4.Range("F12").Value = 372
$ws4.Range("F13").Value = 1094
$ws4.Range("F15").Value = 12410
$ws4.Range("F16").Value = 13151
$ws4.Range("F17").Value = 34
$ws4.Range("F21").Value = 244
$ws4.Range("F24").Value = 130
